$d = $word.ActiveDocument

# 1. JavaScript Skillset line: add "jest, nodemon, socket.io, " before "etc."
$r1 = $d.Content
$ok1 = $r1.Find.Execute(
    "jsonwebtoken, etc.", $true, $false, $false, $false, $false, $true, 1, $false,
    "jsonwebtoken, jest, nodemon, socket.io, etc.", 2)
Write-Output "Step1 (JS skillset): $ok1"

# 2. "Backend Development" -> "Back-End Development"
$r2 = $d.Content
$ok2 = $r2.Find.Execute(
    "Backend Development and General-Purpose Programming", $true, $false, $false, $false, $false, $true, 1, $false,
    "Back-End Development and General-Purpose Programming", 2)
Write-Output "Step2 (Back-End): $ok2"

# 3. "gui/software programs" -> "GUI/software programs"
$r3 = $d.Content
$ok3 = $r3.Find.Execute(
    "gui/software programs", $true, $false, $false, $false, $false, $true, 1, $false,
    "GUI/software programs", 2)
Write-Output "Step3 (GUI): $ok3"

# 4. Append ", git and GitHub work-flow" to the Bash skillset line
$r4 = $d.Content
$ok4 = $r4.Find.Execute(
    "Skillset: Bash shell scripting, network and system management, task automation", $true, $false, $false, $false, $false, $true, 1, $false,
    "Skillset: Bash shell scripting, network and system management, task automation, git and GitHub work-flow", 2)
Write-Output "Step4 (Bash skillset): $ok4"

# 5. Merge "May 2019-" and "Mar 2020" runs into a single run "May 2019-Mar 2020"
$r5 = $d.Content
$ok5 = $r5.Find.Execute(
    "May 2019-Mar 2020", $true, $false, $false, $false, $false, $true, 1, $false,
    "May 2019-Mar 2020", 2)
Write-Output "Step5 (dates merge): $ok5"

# 6. "Great at multi-tasking, ..." -> split into 3 runs with "multitasking" (no hyphen)
$r6 = $d.Content
$ok6 = $r6.Find.Execute(
    "multi-tasking", $true, $false, $false, $false, $false, $true, 1, $false,
    "multitasking", 2)
Write-Output "Step6a (multitasking text): $ok6"

$r6b = $d.Content
$found6b = $r6b.Find.Execute("multitasking", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Output "Step6b (locate multitasking): $found6b Start=$($r6b.Start) End=$($r6b.End)"

$r6c = $d.Range($r6b.Start, $r6b.End)
# Nudge the font size away and back to force Word to keep this span as its
# own run (distinct from the surrounding text) instead of re-merging it.
$r6c.Font.Size = 13
$r6c.Font.Size = 11
Write-Output "Step6c (split run): done"
